$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while forcing it to stay
# plain text (matches source data, which stores Price/Volume as text).
# Writing directly via .Value lets Excel auto-convert digit-dot strings into
# numbers, so instead we build the literal as a text formula in a scratch
# cell (Z1), copy it, and paste-special the *value* into the target cell.
function Set-TextValue($cellRef, $text) {
    $ws.Range("Z1").Formula = '="' + $text + '"'
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("Z1").Clear()
}

Set-TextValue "D2" "42.008.47"
$ws.Range("E2").Value = "  -2.71%  "
Set-TextValue "D3" "2.240.90"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue "D5" "247.23"
$ws.Range("E5").Value = "  -1.88%  "
Set-TextValue "D6" "0.623"
$ws.Range("E6").Value = "  -2.50%  "
Set-TextValue "D7" "76.67"
$ws.Range("E7").Value = "  +4.29%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue "D9" "0.630"
$ws.Range("E9").Value = "  -2.23%  "
Set-TextValue "D10" "40.84"
$ws.Range("E10").Value = "  +4.53%  "
Set-TextValue "D11" "0.0955"
$ws.Range("E11").Value = "  -2.95%  "
Set-TextValue "D12" "7.18"
$ws.Range("E12").Value = "  -3.18%  "
Set-TextValue "D14" "2.574.28"
$ws.Range("E14").Value = "  -2.38%  "
Set-TextValue "D15" "14.87"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("E16").Value = "  -1.99%  "
Set-TextValue "D17" "2.233.23"
$ws.Range("E17").Value = "  -2.69%  "
Set-TextValue "D18" "41.914.64"
$ws.Range("E18").Value = "  -2.65%  "
Set-TextValue "D19" "0.0₃0981"
$ws.Range("E19").Value = "  -2.17%  "
Set-TextValue "D20" "6.14"
$ws.Range("E20").Value = "  -2.72%  "
Set-TextValue "D21" "71.73"
$ws.Range("E21").Value = "  -0.94%  "
Set-TextValue "D22" "2.32"
$ws.Range("E22").Value = "  +3.78%  "
Set-TextValue "D23" "231.64"
$ws.Range("E23").Value = "  -1.37%  "
Set-TextValue "D24" "11.52"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -5.78%  "
Set-TextValue "D27" "2.30"
$ws.Range("E27").Value = "  -5.02%  "
Set-TextValue "D28" "7.20"
$ws.Range("E28").Value = "  +11.69%  "
$ws.Range("E29").Value = "  -1.51%  "
Set-TextValue "D30" "168.61"
$ws.Range("E30").Value = "  +0.93%  "
Set-TextValue "D31" "20.57"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.0836"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D33" "32.83"
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("E34").Value = "  -5.41%  "
Set-TextValue "D35" "0.126"
$ws.Range("E35").Value = "  -0.44%  "
Set-TextValue "D36" "4.54"
$ws.Range("E36").Value = "  -1.10%  "
Set-TextValue "D37" "4.88"
$ws.Range("E37").Value = "  +2.36%  "
Set-TextValue "D38" "0.0301"
$ws.Range("E38").Value = "  -2.45%  "
Set-TextValue "D39" "14.03"
$ws.Range("E39").Value = "  -3.38%  "
Set-TextValue "D40" "5.92"
$ws.Range("E40").Value = "  -0.78%  "
Set-TextValue "D41" "2.19"
$ws.Range("E41").Value = "  -6.79%  "
Set-TextValue "D42" "112.99"
$ws.Range("E42").Value = "  +14.26%  "
Set-TextValue "D43" "0.203"
$ws.Range("E43").Value = "  -7.96%  "
Set-TextValue "D44" "61.15"
$ws.Range("E44").Value = "  -1.35%  "
Set-TextValue "D45" "8.72"
$ws.Range("E45").Value = "  -4.21%  "
Set-TextValue "D46" "0.101"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("E47").Value = "  -0.41%  "
Set-TextValue "D48" "1.14"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D49" "1.17"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D50" "4.30"
$ws.Range("E50").Value = "  -12.93%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D51" "2.26"
$ws.Range("E51").Value = "  -2.24%  "
